# Generate Report for Handback
# Updates the "Latest HO Xliff Generate Date", "Correspond Handoff Datetime"
# and "Correspond Handback DateTime" timestamp cells that get refreshed
# each time the handback status report is regenerated.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Overview sheet: "Latest HO Xliff Generate Date" for the first data row.
$wsOverview.Range("G2").Value = "2016-09-02 05:11:14"

# zh-cn sheet: "Correspond Handoff Datetime" / "Correspond Handback DateTime"
# for the first data row.
$wsZhCn.Range("H2").Value = "2016-09-02 05:11:07"
$wsZhCn.Range("K2").Value = "2016-09-02 05:11:34"

# de-de sheet: same two columns for its first data row.
$wsDeDe.Range("H2").Value = "2016-09-02 05:11:14"
$wsDeDe.Range("K2").Value = "2016-09-02 05:11:42"
